# Apply the three content edits described by the diff:
#  1. Update the letter date from September 19 to September 21, 2025.
#  2. Split the mailing address line into a street line and a city/state/zip
#     line (new paragraph inserted).
#  3. Remove the now-redundant blank NoSpacing paragraph that used to sit
#     right after the "...Board of Directors" line.

$d = $word.ActiveDocument

# 1) Fix the date -----------------------------------------------------
$d.Content.Find.Execute(
    "September 19, 2025", $false, $false, $false, $false, $false,
    $true, 1, $false, "September 21, 2025", 2
) | Out-Null

# 2) Split "999 Story Road, San Jose CA 95122" into two paragraphs ----
$addr = $d.Content
$addr.Find.Execute(
    "999 Story Road, San Jose CA 95122", $false, $false, $false, $false,
    $false, $true, 1, $false, "", 0
) | Out-Null

if ($addr.Find.Found) {
    # Replace the text of the found range with just the street address.
    $addr.Text = "999 Story Road"

    # Insert a new paragraph right after it, carrying the same paragraph
    # formatting, and put the city/state/zip into it.
    $newPara = $addr.Paragraphs(1).Range.InsertParagraphAfter()
    $cityRange = $d.Range($addr.End, $addr.End)
    $cityRange.InsertAfter("San Jose, CA 95122")
}

# 3) Delete the blank paragraph that follows "Board of Directors" ----
$target = $d.Content.Find
$target.Execute("Board of Directors", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($target.Found) {
    $boardRange = $d.Content
    $boardRange.Find.Execute("Board of Directors", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $boardPara = $boardRange.Paragraphs(1)
    $nextPara = $boardPara.Next()
    if ($nextPara -ne $null -and $nextPara.Range.Text.Trim().Length -eq 0) {
        $nextPara.Range.Delete()
    }
}
